# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to the Lamia_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M17").Value = -1032.0
$ws.Range("H17").Value = 3052.487
$ws.Range("J17").Value = 3632.7188
$ws.Range("N17").Value = -11234.1564
$ws.Range("I17").Value = 400.0
$ws.Range("L17").Value = 10898.1564
$ws.Range("K17").Value = 1200.0
$ws.Range("K19").Value = 749.75
$ws.Range("M19").Value = -574.75
$ws.Range("H19").Value = 816.5
$ws.Range("I19").Value = 749.75
$ws.Range("J29").Value = 5781.3335
$ws.Range("H29").Value = 3383.0
$ws.Range("L29").Value = 17344.0005
$ws.Range("N29").Value = -17906.0005
$ws.Range("M32").Value = -10748.0
$ws.Range("H32").Value = 15691.667
$ws.Range("I32").Value = 11074.0
$ws.Range("J32").Value = 18000.5
$ws.Range("L32").Value = 18000.5
$ws.Range("N32").Value = -18652.5
$ws.Range("K32").Value = 11074.0
$ws.Range("N43").Value = -19139.0
$ws.Range("L43").Value = 19001.0
$ws.Range("H43").Value = 19001.0
$ws.Range("J43").Value = 19001.0
$ws.Range("M62").Value = -2118.5625
$ws.Range("H62").Value = 3969.3635
$ws.Range("I62").Value = 2742.5625
$ws.Range("K62").Value = 2742.5625
$ws.Range("M65").Value = -10592.8125
$ws.Range("H65").Value = 3969.3635
$ws.Range("I65").Value = 2742.5625
$ws.Range("K65").Value = 13712.8125
$ws.Range("K98").Value = 544.5455
$ws.Range("M98").Value = 953.4545
$ws.Range("J98").Value = 8443.6
$ws.Range("H98").Value = 3013.0
$ws.Range("I98").Value = 544.5455
$ws.Range("L98").Value = 8443.6
$ws.Range("N98").Value = -11439.6
$ws.Range("N100").Value = -6122.857
$ws.Range("K100").Value = 2038.4
$ws.Range("M100").Value = -1497.4
$ws.Range("J100").Value = 5040.857
$ws.Range("H100").Value = 3789.8333
$ws.Range("I100").Value = 2038.4
$ws.Range("L100").Value = 5040.857
$ws.Range("K112").Value = 1761.0
$ws.Range("M112").Value = -653.0
$ws.Range("H112").Value = 6081.3335
$ws.Range("I112").Value = 587.0
$ws.Range("L113").Value = 7751.5
$ws.Range("N113").Value = -14259.5
$ws.Range("J113").Value = 7751.5
$ws.Range("H113").Value = 6750.625
$ws.Range("K116").Value = 6933.0
$ws.Range("M116").Value = -3491.0
$ws.Range("H116").Value = 7127.68
$ws.Range("I116").Value = 6933.0
$ws.Range("M122").Value = 816.3635000000002
$ws.Range("J122").Value = 8443.6
$ws.Range("H122").Value = 3013.0
$ws.Range("I122").Value = 544.5455
$ws.Range("L122").Value = 25330.8
$ws.Range("N122").Value = -30230.8
$ws.Range("K122").Value = 1633.6365
$ws.Range("K132").Value = 5709.0
$ws.Range("M132").Value = -3179.0
$ws.Range("H132").Value = 1903.0
$ws.Range("I132").Value = 1903.0
$ws.Range("H137").Value = 3874.5557
$ws.Range("I137").Value = 0.0
$ws.Range("J137").Value = 3874.5557
$ws.Range("L137").Value = 11623.6671
$ws.Range("N137").Value = -16723.6671
$ws.Range("K137").Value = 0.0
$ws.Range("H138").Value = 3209.5312
$ws.Range("J138").Value = 3200.85
$ws.Range("N138").Value = -19882.55
$ws.Range("I138").Value = 3224.0
$ws.Range("L138").Value = 9602.55
$ws.Range("K138").Value = 9672.0
$ws.Range("M138").Value = -4532.0
$ws.Range("H141").Value = 8364.9
$ws.Range("I141").Value = 7824.5
$ws.Range("L141").Value = 25500.0
$ws.Range("N141").Value = -35860.0
$ws.Range("K141").Value = 23473.5
$ws.Range("M141").Value = -18293.5
$ws.Range("J141").Value = 8500.0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 456.46155
$ws.Range("J5").Value = 678.125
$ws.Range("L5").Value = 678.125
$ws.Range("N5").Value = -902.125
$ws.Range("H32").Value = 4620.095
$ws.Range("J32").Value = 10045.4
$ws.Range("L32").Value = 10045.4
$ws.Range("N32").Value = -10619.4
$ws.Range("M45").Value = -142858223.0
$ws.Range("H45").Value = 71433470.0
$ws.Range("I45").Value = 142858600.0
$ws.Range("K45").Value = 142858600.0
$ws.Range("J61").Value = 4313.7144
$ws.Range("H61").Value = 3037.1353
$ws.Range("I61").Value = 2739.2666
$ws.Range("L61").Value = 4313.7144
$ws.Range("N61").Value = -4737.7144
$ws.Range("K61").Value = 2739.2666
$ws.Range("M61").Value = -2527.2666
$ws.Range("K74").Value = 22224366.0
$ws.Range("M74").Value = -22223492.0
$ws.Range("H74").Value = 18520552.0
$ws.Range("J74").Value = 1478.0
$ws.Range("N74").Value = -3226.0
$ws.Range("I74").Value = 22224366.0
$ws.Range("L74").Value = 1478.0
$ws.Range("I77").Value = 22224366.0
$ws.Range("J77").Value = 1478.0
$ws.Range("L77").Value = 7390.0
$ws.Range("N77").Value = -16126.0
$ws.Range("K77").Value = 111121830.0
$ws.Range("M77").Value = -111117462.0
$ws.Range("H77").Value = 18520552.0
$ws.Range("K110").Value = 677.5
$ws.Range("M110").Value = 1367.5
$ws.Range("H110").Value = 1487.591
$ws.Range("I110").Value = 677.5
$ws.Range("H118").Value = 0.0
$ws.Range("J118").Value = 0.0
$ws.Range("L118").Value = 0.0
$ws.Range("M122").Value = -2827.800099999999
$ws.Range("H122").Value = 2274.3125
$ws.Range("I122").Value = 1759.2667
$ws.Range("K122").Value = 5277.800099999999
$ws.Range("K132").Value = 8045.646900000001
$ws.Range("M132").Value = -5515.646900000001
$ws.Range("J132").Value = 13592.667
$ws.Range("H132").Value = 4318.5
$ws.Range("I132").Value = 2681.8823
$ws.Range("L132").Value = 40778.001
$ws.Range("N132").Value = -45838.001
$ws.Range("K136").Value = 8217.7998
$ws.Range("M136").Value = -5667.799800000001
$ws.Range("H136").Value = 3037.1353
$ws.Range("I136").Value = 2739.2666
$ws.Range("J136").Value = 4313.7144
$ws.Range("L136").Value = 12941.1432
$ws.Range("N136").Value = -18041.1432
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N4").Value = -908.125
$ws.Range("J4").Value = 678.125
$ws.Range("H4").Value = 456.46155
$ws.Range("L4").Value = 678.125
$ws.Range("M22").Value = -151.75
$ws.Range("H22").Value = 659.8
$ws.Range("I22").Value = 324.75
$ws.Range("K22").Value = 324.75
$ws.Range("K94").Value = 867.5
$ws.Range("M94").Value = -416.5
$ws.Range("H94").Value = 955.8261
$ws.Range("J94").Value = 2899.0
$ws.Range("N94").Value = -3801.0
$ws.Range("I94").Value = 867.5
$ws.Range("L94").Value = 2899.0
$ws.Range("K99").Value = 1316.8
$ws.Range("M99").Value = 181.2
$ws.Range("H99").Value = 1379.0
$ws.Range("I99").Value = 1316.8
$ws.Range("M105").Value = -27232.0
$ws.Range("J105").Value = 11999.667
$ws.Range("H105").Value = 19989.941
$ws.Range("I105").Value = 28979.0
$ws.Range("L105").Value = 11999.667
$ws.Range("N105").Value = -15493.667
$ws.Range("K105").Value = 28979.0
$ws.Range("M107").Value = 710.8462
$ws.Range("H107").Value = 1209.1538
$ws.Range("I107").Value = 1209.1538
$ws.Range("K107").Value = 1209.1538
$ws.Range("M134").Value = -3770.499899999999
$ws.Range("H134").Value = 4702.875
$ws.Range("I134").Value = 2101.8333
$ws.Range("J134").Value = 12506.0
$ws.Range("L134").Value = 37518.0
$ws.Range("N134").Value = -42588.0
$ws.Range("K134").Value = 6305.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J7").Value = 376.27274
$ws.Range("H7").Value = 721.7273
$ws.Range("I7").Value = 1067.1818
$ws.Range("L7").Value = 376.27274
$ws.Range("N7").Value = -602.27274
$ws.Range("K7").Value = 1067.1818
$ws.Range("M7").Value = -954.1818000000001
$ws.Range("M22").Value = 181.09091
$ws.Range("H22").Value = 667.2
$ws.Range("I22").Value = 168.90909
$ws.Range("J22").Value = 2037.5
$ws.Range("N22").Value = -2737.5
$ws.Range("L22").Value = 2037.5
$ws.Range("K22").Value = 168.90909
$ws.Range("H31").Value = 29806.25
$ws.Range("J31").Value = 92086.0
$ws.Range("N31").Value = -92676.0
$ws.Range("I31").Value = 3114.9285
$ws.Range("L31").Value = 92086.0
$ws.Range("K31").Value = 3114.9285
$ws.Range("M31").Value = -2819.9285
$ws.Range("M34").Value = -2912.9285
$ws.Range("H34").Value = 29806.25
$ws.Range("I34").Value = 3114.9285
$ws.Range("J34").Value = 92086.0
$ws.Range("N34").Value = -92490.0
$ws.Range("L34").Value = 92086.0
$ws.Range("K34").Value = 3114.9285
$ws.Range("N43").Value = -8368.0
$ws.Range("L43").Value = 8000.0
$ws.Range("H43").Value = 8000.0
$ws.Range("J43").Value = 8000.0
$ws.Range("M45").Value = -1574.0
$ws.Range("H45").Value = 4883.5
$ws.Range("I45").Value = 2167.0
$ws.Range("K45").Value = 2167.0
$ws.Range("H50").Value = 40883.332
$ws.Range("J50").Value = 41135.715
$ws.Range("L50").Value = 41135.715
$ws.Range("N50").Value = -42385.715
$ws.Range("M55").Value = -8077.25
$ws.Range("H55").Value = 12617.556
$ws.Range("I55").Value = 8392.25
$ws.Range("K55").Value = 8392.25
$ws.Range("K82").Value = 0.0
$ws.Range("J82").Value = 38999.0
$ws.Range("H82").Value = 38999.0
$ws.Range("I82").Value = 0.0
$ws.Range("L82").Value = 38999.0
$ws.Range("N82").Value = -39721.0
$ws.Range("I85").Value = 0.0
$ws.Range("L85").Value = 38999.0
$ws.Range("N85").Value = -41495.0
$ws.Range("K85").Value = 0.0
$ws.Range("H85").Value = 38999.0
$ws.Range("J85").Value = 38999.0
$ws.Range("K86").Value = 5374.625
$ws.Range("M86").Value = -4251.625
$ws.Range("H86").Value = 6975.3335
$ws.Range("I86").Value = 5374.625
$ws.Range("K89").Value = 26873.125
$ws.Range("M89").Value = -21257.125
$ws.Range("H89").Value = 6975.3335
$ws.Range("I89").Value = 5374.625
$ws.Range("K99").Value = 3383.25
$ws.Range("M99").Value = -1885.25
$ws.Range("J99").Value = 3400.0
$ws.Range("H99").Value = 3386.6
$ws.Range("I99").Value = 3383.25
$ws.Range("L99").Value = 3400.0
$ws.Range("N99").Value = -6396.0
$ws.Range("H101").Value = 8000.0
$ws.Range("J101").Value = 8000.0
$ws.Range("L101").Value = 8000.0
$ws.Range("N101").Value = -14490.0
$ws.Range("M105").Value = -3419.0
$ws.Range("H105").Value = 6000.778
$ws.Range("I105").Value = 5166.0
$ws.Range("K105").Value = 5166.0
$ws.Range("M107").Value = 1111.05
$ws.Range("H107").Value = 1050.037
$ws.Range("I107").Value = 808.95
$ws.Range("J107").Value = 1738.8572
$ws.Range("N107").Value = -5578.8572
$ws.Range("L107").Value = 1738.8572
$ws.Range("K107").Value = 808.95
$ws.Range("N126").Value = -15140.0
$ws.Range("K126").Value = 10149.75
$ws.Range("M126").Value = -7679.75
$ws.Range("H126").Value = 3386.6
$ws.Range("I126").Value = 3383.25
$ws.Range("J126").Value = 3400.0
$ws.Range("L126").Value = 10200.0
$ws.Range("K132").Value = 13541.5005
$ws.Range("M132").Value = -11011.5005
$ws.Range("H132").Value = 4690.5933
$ws.Range("I132").Value = 4513.8335
$ws.Range("M134").Value = -8241.0
$ws.Range("H134").Value = 5587.25
$ws.Range("I134").Value = 3592.0
$ws.Range("J134").Value = 8912.667
$ws.Range("L134").Value = 26738.001
$ws.Range("N134").Value = -31808.001
$ws.Range("K134").Value = 10776.0
$ws.Range("H135").Value = 69996.664
$ws.Range("J135").Value = 69996.664
$ws.Range("L135").Value = 69996.664
$ws.Range("N135").Value = -80136.664
$ws.Range("H141").Value = 239993.6
$ws.Range("L141").Value = 239993.6
$ws.Range("N141").Value = -250353.6
$ws.Range("J141").Value = 239993.6
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J86").Value = 2982.6667
$ws.Range("H86").Value = 2148.4443
$ws.Range("L86").Value = 8948.000100000001
$ws.Range("N86").Value = -11320.0001
$ws.Range("N89").Value = -38700.0003
$ws.Range("H89").Value = 2148.4443
$ws.Range("J89").Value = 2982.6667
$ws.Range("L89").Value = 26844.0003
$ws.Range("H92").Value = 2733.5557
$ws.Range("L92").Value = 19801.9995
$ws.Range("N92").Value = -22297.9995
$ws.Range("J92").Value = 6600.6665
$ws.Range("J104").Value = 6174.5
$ws.Range("H104").Value = 6174.5
$ws.Range("L104").Value = 18523.5
$ws.Range("N104").Value = -23765.5
$ws.Range("H109").Value = 3608.1
$ws.Range("I109").Value = 3731.3333
$ws.Range("J109").Value = 2499.0
$ws.Range("L109").Value = 7497.0
$ws.Range("N109").Value = -9577.0
$ws.Range("K109").Value = 11193.9999
$ws.Range("M109").Value = -10153.9999
$ws.Range("K132").Value = 32254.4997
$ws.Range("M132").Value = -29724.4997
$ws.Range("J132").Value = 4403.9165
$ws.Range("H132").Value = 4130.5557
$ws.Range("I132").Value = 3583.8333
$ws.Range("L132").Value = 39635.2485
$ws.Range("N132").Value = -44695.2485
$ws.Range("K136").Value = 4590.0
$ws.Range("M136").Value = 510.0
$ws.Range("H136").Value = 4515.0
$ws.Range("I136").Value = 1530.0
$ws.Range("J136").Value = 7500.0
$ws.Range("L136").Value = 22500.0
$ws.Range("N136").Value = -32700.0
$ws.Range("H138").Value = 25513.334
$ws.Range("J138").Value = 16250.0
$ws.Range("N138").Value = -59030.0
$ws.Range("I138").Value = 30145.0
$ws.Range("L138").Value = 48750.0
$ws.Range("K138").Value = 90435.0
$ws.Range("M138").Value = -85295.0
$ws.Range("H139").Value = 7446.25
$ws.Range("I139").Value = 2205.0
$ws.Range("J139").Value = 11190.0
$ws.Range("N139").Value = -43850.0
$ws.Range("L139").Value = 33570.0
$ws.Range("K139").Value = 6615.0
$ws.Range("M139").Value = -1475.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K2").Value = 20.25
$ws.Range("M2").Value = 92.75
$ws.Range("H2").Value = 941.6429
$ws.Range("I2").Value = 20.25
$ws.Range("H26").Value = 27449.5
$ws.Range("I26").Value = 0.0
$ws.Range("K26").Value = 0.0
$ws.Range("M34").Value = -29732.0
$ws.Range("H34").Value = 30000.0
$ws.Range("I34").Value = 30000.0
$ws.Range("K34").Value = 30000.0
$ws.Range("H50").Value = 27449.5
$ws.Range("I50").Value = 0.0
$ws.Range("K50").Value = 0.0
$ws.Range("M76").Value = -29685.0
$ws.Range("H76").Value = 30000.0
$ws.Range("I76").Value = 30000.0
$ws.Range("K76").Value = 30000.0
$ws.Range("H79").Value = 30000.0
$ws.Range("I79").Value = 30000.0
$ws.Range("K79").Value = 30000.0
$ws.Range("M79").Value = -28908.0
$ws.Range("H80").Value = 10593.333
$ws.Range("J80").Value = 10593.333
$ws.Range("L80").Value = 10593.333
$ws.Range("N80").Value = -12589.333
$ws.Range("H83").Value = 10593.333
$ws.Range("J83").Value = 10593.333
$ws.Range("L83").Value = 52966.665
$ws.Range("N83").Value = -62950.665
$ws.Range("N100").Value = -22164.0
$ws.Range("J100").Value = 20000.0
$ws.Range("H100").Value = 20000.0
$ws.Range("L100").Value = 20000.0
$ws.Range("N102").Value = -22257.0
$ws.Range("K102").Value = 1778.0714
$ws.Range("M102").Value = -156.0714
$ws.Range("H102").Value = 2927.0667
$ws.Range("I102").Value = 1778.0714
$ws.Range("J102").Value = 19013.0
$ws.Range("L102").Value = 19013.0
$ws.Range("K132").Value = 437270.16
$ws.Range("M132").Value = -434740.16
$ws.Range("H132").Value = 105551.1
$ws.Range("I132").Value = 145756.72
$ws.Range("M26").ClearContents()
$ws.Range("M50").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J46").Value = 4332.6665
$ws.Range("H46").Value = 3772.2727
$ws.Range("L46").Value = 4332.6665
$ws.Range("N46").Value = -4708.6665
$ws.Range("M55").Value = -2501164.5
$ws.Range("H55").Value = 1516967.4
$ws.Range("I55").Value = 2501337.5
$ws.Range("J55").Value = 2551.923
$ws.Range("L55").Value = 2551.923
$ws.Range("N55").Value = -2897.923
$ws.Range("K55").Value = 2501337.5
$ws.Range("J61").Value = 10150.25
$ws.Range("H61").Value = 2873.2104
$ws.Range("L61").Value = 10150.25
$ws.Range("N61").Value = -10554.25
$ws.Range("K74").Value = 45000.0
$ws.Range("M74").Value = -44002.0
$ws.Range("H74").Value = 45000.0
$ws.Range("I74").Value = 45000.0
$ws.Range("I77").Value = 45000.0
$ws.Range("K77").Value = 135000.0
$ws.Range("M77").Value = -130008.0
$ws.Range("H77").Value = 45000.0
$ws.Range("L113").Value = 10150.25
$ws.Range("N113").Value = -14490.25
$ws.Range("J113").Value = 10150.25
$ws.Range("H113").Value = 2873.2104
$ws.Range("M122").Value = -21159.25
$ws.Range("J122").Value = 15002.5
$ws.Range("H122").Value = 9296.3
$ws.Range("I122").Value = 7869.75
$ws.Range("L122").Value = 45007.5
$ws.Range("N122").Value = -49907.5
$ws.Range("K122").Value = 23609.25
$ws.Range("K132").Value = 8097.0
$ws.Range("M132").Value = -5567.0
$ws.Range("H132").Value = 6299.9287
$ws.Range("I132").Value = 2699.0
$ws.Range("K136").Value = 8831.7276
$ws.Range("M136").Value = -6281.7276
$ws.Range("H136").Value = 6783.0
$ws.Range("I136").Value = 2943.9092
$ws.Range("J136").Value = 11006.0
$ws.Range("L136").Value = 33018.0
$ws.Range("N136").Value = -38118.0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K70").Value = 0.0
$ws.Range("H70").Value = 0.0
$ws.Range("I70").Value = 0.0
$ws.Range("K73").Value = 0.0
$ws.Range("H73").Value = 0.0
$ws.Range("I73").Value = 0.0
$ws.Range("K82").Value = 50000.0
$ws.Range("M82").Value = -49617.0
$ws.Range("H82").Value = 50000.0
$ws.Range("I82").Value = 50000.0
$ws.Range("I85").Value = 50000.0
$ws.Range("K85").Value = 50000.0
$ws.Range("M85").Value = -48674.0
$ws.Range("H85").Value = 50000.0
$ws.Range("K96").Value = 2795.8
$ws.Range("M96").Value = -1422.8
$ws.Range("J96").Value = 9999.667
$ws.Range("H96").Value = 5497.25
$ws.Range("I96").Value = 2795.8
$ws.Range("L96").Value = 9999.667
$ws.Range("N96").Value = -12745.667
$ws.Range("M122").Value = -3266.799999999999
$ws.Range("H122").Value = 7570.567
$ws.Range("I122").Value = 1905.6
$ws.Range("K122").Value = 5716.799999999999
$ws.Range("N123").Value = -38800.0
$ws.Range("J123").Value = 29000.0
$ws.Range("H123").Value = 29000.0
$ws.Range("L123").Value = 29000.0
$ws.Range("K126").Value = 3543.0
$ws.Range("M126").Value = -1073.0
$ws.Range("H126").Value = 2533.375
$ws.Range("I126").Value = 1181.0
$ws.Range("K132").Value = 15655.0005
$ws.Range("M132").Value = -13125.0005
$ws.Range("H132").Value = 5378.75
$ws.Range("I132").Value = 5218.3335
$ws.Range("K136").Value = 6714.545999999999
$ws.Range("M136").Value = -4164.545999999999
$ws.Range("H136").Value = 3635.4167
$ws.Range("I136").Value = 2238.182
$ws.Range("H138").Value = 0.0
$ws.Range("J138").Value = 0.0
$ws.Range("L138").Value = 0.0
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N138").ClearContents()

